$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1894736842105263
$ws.Range("C2").Value = 0.5543859649122806
$ws.Range("J2").Value = 0.02456140350877193
$ws.Range("P2").Value = 0.1403508771929824
$ws.Range("S2").Value = 0.0912280701754386
$ws.Range("B3").Value = 0.02941176470588235
$ws.Range("C3").Value = 0.05882352941176471
$ws.Range("J3").Value = 0.02352941176470588
$ws.Range("O3").Value = 0.005882352941176471
$ws.Range("P3").Value = 0.7058823529411765
$ws.Range("S3").Value = 0.1764705882352941
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.5294117647058824
$ws.Range("S4").Value = 0.4411764705882353
$ws.Range("P5").Value = 0.8888888888888888
$ws.Range("S5").Value = 0.1111111111111111
$ws.Range("B6").Value = 0.04977375565610859
$ws.Range("F6").Value = 0.04072398190045249
$ws.Range("J6").Value = 0.2398190045248869
$ws.Range("O6").Value = 0.03167420814479638
$ws.Range("Q6").Value = 0.1402714932126697
$ws.Range("R6").Value = 0.07239819004524888
$ws.Range("S6").Value = 0.4253393665158371
$ws.Range("B7").Value = 0.09547738693467336
$ws.Range("D7").Value = 0.01507537688442211
$ws.Range("E7").Value = 0.01005025125628141
$ws.Range("F7").Value = 0.05025125628140704
$ws.Range("J7").Value = 0.1407035175879397
$ws.Range("O7").Value = 0.03517587939698492
$ws.Range("Q7").Value = 0.1557788944723618
$ws.Range("R7").Value = 0.1055276381909548
$ws.Range("S7").Value = 0.3919597989949749
$ws.Range("B8").Value = 0.07586206896551724
$ws.Range("D8").Value = 0.01609195402298851
$ws.Range("F8").Value = 0.05977011494252873
$ws.Range("J8").Value = 0.128735632183908
$ws.Range("O8").Value = 0.02528735632183908
$ws.Range("Q8").Value = 0.1862068965517241
$ws.Range("R8").Value = 0.09425287356321839
$ws.Range("S8").Value = 0.4137931034482759
$ws.Range("B9").Value = 0.1054852320675106
$ws.Range("D9").Value = 0.008438818565400843
$ws.Range("F9").Value = 0.06329113924050633
$ws.Range("J9").Value = 0.1265822784810127
$ws.Range("O9").Value = 0.01687763713080169
$ws.Range("Q9").Value = 0.1687763713080169
$ws.Range("R9").Value = 0.1181434599156118
$ws.Range("S9").Value = 0.3924050632911392
$ws.Range("B10").Value = 0.1048573631457209
$ws.Range("D10").Value = 0.01773323053199691
$ws.Range("E10").Value = 0.006168080185042405
$ws.Range("F10").Value = 0.07324595219737856
$ws.Range("J10").Value = 0.1148804934464148
$ws.Range("O10").Value = 0.01850424055512722
$ws.Range("Q10").Value = 0.2158828064764842
$ws.Range("R10").Value = 0.07478797224363917
$ws.Range("S10").Value = 0.3739398612181958
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("J11").Value = 0.0915032679738562
$ws.Range("K11").Value = 0.1535947712418301
$ws.Range("L11").Value = 0.6143790849673203
$ws.Range("S11").Value = 0.02287581699346405
$ws.Range("G12").Value = 0.7135678391959799
$ws.Range("J12").Value = 0.1959798994974874
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.03517587939698492
$ws.Range("S12").Value = 0.05025125628140704
$ws.Range("G13").Value = 0.7027027027027027
$ws.Range("J13").Value = 0.1891891891891892
$ws.Range("S13").Value = 0.1081081081081081
$ws.Range("F15").Value = 0.02173913043478261
$ws.Range("H15").Value = 0.1695652173913043
$ws.Range("I15").Value = 0.05652173913043478
$ws.Range("J15").Value = 0.3434782608695652
$ws.Range("K15").Value = 0.07391304347826087
$ws.Range("O15").Value = 0.06521739130434782
$ws.Range("S15").Value = 0.2695652173913043
$ws.Range("F16").Value = 0.02209944751381215
$ws.Range("H16").Value = 0.1325966850828729
$ws.Range("I16").Value = 0.143646408839779
$ws.Range("J16").Value = 0.3370165745856354
$ws.Range("K16").Value = 0.08839779005524862
$ws.Range("M16").Value = 0.03314917127071823
$ws.Range("O16").Value = 0.04419889502762431
$ws.Range("S16").Value = 0.1988950276243094
$ws.Range("F17").Value = 0.01694915254237288
$ws.Range("H17").Value = 0.1885593220338983
$ws.Range("I17").Value = 0.1059322033898305
$ws.Range("J17").Value = 0.413135593220339
$ws.Range("K17").Value = 0.08898305084745763
$ws.Range("M17").Value = 0.01694915254237288
$ws.Range("N17").Value = 0.00423728813559322
$ws.Range("O17").Value = 0.04449152542372881
$ws.Range("S17").Value = 0.1207627118644068
$ws.Range("F18").Value = 0.004830917874396135
$ws.Range("H18").Value = 0.1884057971014493
$ws.Range("I18").Value = 0.1449275362318841
$ws.Range("J18").Value = 0.357487922705314
$ws.Range("K18").Value = 0.0821256038647343
$ws.Range("M18").Value = 0.01449275362318841
$ws.Range("O18").Value = 0.06280193236714976
$ws.Range("S18").Value = 0.1449275362318841
$ws.Range("F19").Value = 0.01572551822730522
$ws.Range("H19").Value = 0.1801286633309507
$ws.Range("I19").Value = 0.08506075768406005
$ws.Range("J19").Value = 0.359542530378842
$ws.Range("K19").Value = 0.1150822015725518
$ws.Range("M19").Value = 0.01572551822730522
$ws.Range("N19").Value = 0.002144388849177984
$ws.Range("O19").Value = 0.06790564689063616
$ws.Range("S19").Value = 0.1586847748391708

